# Handle the Female and Male learner comments and editing
#
# Sheet1 had columns D:F as  learners_no | female_learners | male_learners.
# Re-order them so the breakdown comes first and the total follows:
#       D:F -> female_learners | male_learners | learners_no
# and turn the (previously hand-typed) total into a live formula that
# sums the female/male counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Move column E (female_learners) so it lands before column D (learners_no).
$ws.Columns.Item(5).Cut()
$ws.Columns.Item(4).Insert()

# Column F (male_learners) is still in place; move it so it lands before the
# (now 3rd) learners_no column, giving the final D,E,F order of
# female_learners, male_learners, learners_no.
$ws.Columns.Item(6).Cut()
$ws.Columns.Item(5).Insert()

# The new trailing column (F) carries the old learners_no header/width but
# should show the default width used by the female/male columns.
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668

# learners_no (now F2) becomes a formula totalling the two new columns
# instead of a manually entered number.
$ws.Range("F2").Formula = "=SUM(D2+E2)"

# Leave the cursor where the editor finished working.
$ws.Range("F8").Select()
